# Update "想去人数" (want-to-go count) figures in both the "展览" and
# "全部类型" worksheets, which carry duplicate copies of the same data.
#
#   Row 2 (F2): 529 -> 532
#   Row 4 (F4): 168 -> 169
#   Row 5 (F5):   1 -> 2
#   Row 7 (F7): 739 -> 746

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2 = 532
    4 = 169
    5 = 2
    7 = 746
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
